# Femacal de La Calera - Zanahoria: weekly update
# Insert two new price records as rows 177 and 178 (pushing the existing
# rows 177-224 down to 179-226), matching the new dimension A1:R226.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 177 (so the old row 177 ends up at row 179).
$ws.Rows.Item(177).Insert()
$ws.Rows.Item(177).Insert()

# New row 177
$ws.Cells.Item(177, 1).Value = 3
$ws.Cells.Item(177, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(177, 3).Value = "Coquimbo"
$ws.Cells.Item(177, 4).Value = 44508
$ws.Cells.Item(177, 5).Value = 5
$ws.Cells.Item(177, 6).Value = 100114013
$ws.Cells.Item(177, 7).Value = "Zanahoria"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 310
$ws.Cells.Item(177, 11).Value = 6500
$ws.Cells.Item(177, 12).Value = 7000
$ws.Cells.Item(177, 13).Value = 6758
$ws.Cells.Item(177, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(177, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(177, 16).Value = 338
$ws.Cells.Item(177, 17).Value = 20
$ws.Cells.Item(177, 18).Value = "Hortaliza"

# New row 178
$ws.Cells.Item(178, 1).Value = 3
$ws.Cells.Item(178, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(178, 3).Value = "Coquimbo"
$ws.Cells.Item(178, 4).Value = 44508
$ws.Cells.Item(178, 5).Value = 5
$ws.Cells.Item(178, 6).Value = 100114013
$ws.Cells.Item(178, 7).Value = "Zanahoria"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Segunda"
$ws.Cells.Item(178, 10).Value = 160
$ws.Cells.Item(178, 11).Value = 4500
$ws.Cells.Item(178, 12).Value = 4500
$ws.Cells.Item(178, 13).Value = 4500
$ws.Cells.Item(178, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(178, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(178, 16).Value = 225
$ws.Cells.Item(178, 17).Value = 20
$ws.Cells.Item(178, 18).Value = "Hortaliza"

# Make sure the date cells keep the original date number format/style (s="2"),
# same as every other "Fecha" column cell in the sheet.
$ws.Cells.Item(177, 4).NumberFormat = $ws.Cells.Item(179, 4).NumberFormat
$ws.Cells.Item(178, 4).NumberFormat = $ws.Cells.Item(179, 4).NumberFormat
